$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.455.21"

$ws.Range("D3").Value = "2.417.53"
$ws.Range("E3").Value = "  +2.21%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.512"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.07%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  +4.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.20"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.03%  "

$ws.Range("E11").Value = "  +1.61%  "

$ws.Range("E12").Value = "  +4.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.123"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.00%  "

$ws.Range("E14").Value = "  +2.84%  "

$ws.Range("D15").Value = "2.793.89"
$ws.Range("E15").Value = "  +2.35%  "

$ws.Range("D16").Value = "2.428.76"
$ws.Range("E16").Value = "  +3.15%  "

$ws.Range("E17").Value = "  +4.79%  "

$ws.Range("D18").Value = "44.303.58"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.00%  "

$ws.Range("E20").Value = "  +1.63%  "

$ws.Range("D21").Value = "0.0₃0920"
$ws.Range("E21").Value = "  +3.79%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.34%  "

$ws.Range("E25").Value = "  +1.62%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.46%  "

$ws.Range("E28").Value = "  -4.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.54%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "48.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.79%  "

$ws.Range("E32").Value = "  +18.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.34%  "

$ws.Range("E34").Value = "  +2.61%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0773"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.97%  "

$ws.Range("E36").Value = "  +0.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.87"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.48%  "

$ws.Range("E38").Value = "  +2.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.37%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "121.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.85%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.73%  "

$ws.Range("E42").Value = "  +1.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.47%  "

$ws.Range("E44").Value = "  +3.42%  "

$ws.Range("D45").Value = "1.944.39"
$ws.Range("E45").Value = "  +0.56%  "

$ws.Range("E46").Value = "  +1.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.98%  "

$ws.Range("E48").Value = "  +3.43%  "

$ws.Range("E49").Value = "  +9.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "55.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.79%  "
